$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("I5").Value = 2.8
$ws.Range("V5").Value = 9.5
$ws.Range("W5").Value = 23
$ws.Range("AG5").Value = 10
$ws.Range("AJ5").Value = 26

# Row 6
$ws.Range("I6").Value = 6.6
$ws.Range("AH6").Value = 175

# Row 7
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3.55
$ws.Range("I7").Value = 4.6
$ws.Range("L7").Value = 1.29
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 1.85
$ws.Range("O7").Value = 1.75
$ws.Range("P7").Value = 1.4
$ws.Range("Q7").Value = 2.52
$ws.Range("R7").Value = 1.82
$ws.Range("S7").Value = 1.78
$ws.Range("T7").Value = 6.6
$ws.Range("U7").Value = 7.7
$ws.Range("V7").Value = 8.25
$ws.Range("W7").Value = 13
$ws.Range("X7").Value = 14
$ws.Range("Y7").Value = 28
$ws.Range("Z7").Value = 9.5
$ws.Range("AA7").Value = 7
$ws.Range("AB7").Value = 16.5
$ws.Range("AC7").Value = 80
$ws.Range("AD7").Value = 700
$ws.Range("AE7").Value = 12
$ws.Range("AF7").Value = 26
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 80
$ws.Range("AJ7").Value = 50

# Row 8
$ws.Range("G8").Value = 1.36
$ws.Range("H8").Value = 4.35
$ws.Range("I8").Value = 8
$ws.Range("L8").Value = 1.22
$ws.Range("M8").Value = 3.45
$ws.Range("N8").Value = 1.65
$ws.Range("O8").Value = 1.98
$ws.Range("P8").Value = 1.39
$ws.Range("Q8").Value = 2.55
$ws.Range("R8").Value = 1.98
$ws.Range("S8").Value = 1.65
$ws.Range("T8").Value = 6.7
$ws.Range("U8").Value = 6.3
$ws.Range("V8").Value = 8.25
$ws.Range("W8").Value = 8.5
$ws.Range("X8").Value = 11.5
$ws.Range("Y8").Value = 28
$ws.Range("Z8").Value = 11.75
$ws.Range("AA8").Value = 8.75
$ws.Range("AB8").Value = 20
$ws.Range("AC8").Value = 100
$ws.Range("AD8").Value = 800
$ws.Range("AE8").Value = 21
$ws.Range("AF8").Value = 55
$ws.Range("AI8").Value = 90
$ws.Range("AJ8").Value = 75

# Row 9
$ws.Range("G9").Value = 1.87
$ws.Range("H9").Value = 3.3
$ws.Range("L9").Value = 1.37
$ws.Range("M9").Value = 2.65
$ws.Range("N9").Value = 2.07
$ws.Range("O9").Value = 1.6
$ws.Range("P9").Value = 1.47
$ws.Range("Q9").Value = 2.35
$ws.Range("R9").Value = 1.91
$ws.Range("T9").Value = 6.1
$ws.Range("U9").Value = 8
$ws.Range("W9").Value = 15.5
$ws.Range("X9").Value = 16
$ws.Range("Y9").Value = 32
$ws.Range("Z9").Value = 8
$ws.Range("AA9").Value = 6.4
$ws.Range("AD9").Value = 900
$ws.Range("AE9").Value = 9.75
$ws.Range("AF9").Value = 21
$ws.Range("AJ9").Value = 55

# Row 10
$ws.Range("H10").Value = 3.45
$ws.Range("I10").Value = 2.15
$ws.Range("L10").Value = 1.29
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 1.85
$ws.Range("O10").Value = 1.75
$ws.Range("P10").Value = 1.39
$ws.Range("Q10").Value = 2.57
$ws.Range("R10").Value = 1.72
$ws.Range("S10").Value = 1.88
$ws.Range("T10").Value = 9.5
$ws.Range("U10").Value = 15.5
$ws.Range("V10").Value = 11
$ws.Range("X10").Value = 26
$ws.Range("Y10").Value = 35
$ws.Range("Z10").Value = 10.25
$ws.Range("AA10").Value = 6.7
$ws.Range("AB10").Value = 14.5
$ws.Range("AC10").Value = 70
$ws.Range("AD10").Value = 500
$ws.Range("AE10").Value = 7.7
$ws.Range("AF10").Value = 10.25
$ws.Range("AG10").Value = 9
$ws.Range("AH10").Value = 20
$ws.Range("AI10").Value = 17.5
$ws.Range("AJ10").Value = 29

# Row 13
$ws.Range("G13").Value = 2.22
$ws.Range("H13").Value = 2.62
$ws.Range("J13").Value = 1.19
$ws.Range("K13").Value = 4.15
$ws.Range("L13").Value = 1.8
$ws.Range("M13").Value = 1.91
$ws.Range("N13").Value = 3.3
$ws.Range("O13").Value = 1.28
$ws.Range("P13").Value = 1.75
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 2.57
$ws.Range("S13").Value = 1.44
$ws.Range("V13").Value = 10.75
$ws.Range("W13").Value = 22
$ws.Range("Z13").Value = 4.15
$ws.Range("AB13").Value = 26
$ws.Range("AE13").Value = 6.9
$ws.Range("AG13").Value = 16.5
$ws.Range("AH13").Value = 75
$ws.Range("AI13").Value = 65
$ws.Range("AJ13").Value = 100

# Row 19
$ws.Range("G19").Value = 5.1
$ws.Range("I19").Value = 1.57
$ws.Range("L19").Value = 1.2
$ws.Range("O19").Value = 2.18
$ws.Range("P19").Value = 1.33
$ws.Range("Q19").Value = 3.05
$ws.Range("T19").Value = 16.5
$ws.Range("U19").Value = 32
$ws.Range("V19").Value = 16
$ws.Range("W19").Value = 90
$ws.Range("X19").Value = 45
$ws.Range("Y19").Value = 45
$ws.Range("AB19").Value = 14.5
$ws.Range("AC19").Value = 55
$ws.Range("AE19").Value = 8.25
$ws.Range("AF19").Value = 8.25
$ws.Range("AH19").Value = 12
$ws.Range("AI19").Value = 11.5
$ws.Range("AJ19").Value = 21

# Row 25
$ws.Range("J25").Value = 1.07
$ws.Range("K25").Value = 9
$ws.Range("L25").Value = 1.33
$ws.Range("M25").Value = 3.25
$ws.Range("N25").Value = 2.08
$ws.Range("O25").Value = 1.73
$ws.Range("R25").Value = 1.8
$ws.Range("S25").Value = 1.95
$ws.Range("X25").Value = 21
$ws.Range("Z25").Value = 9
$ws.Range("AI25").Value = 26

# Row 28
$ws.Range("N28").Value = 1.65
$ws.Range("O28").Value = 2.2
$ws.Range("W28").Value = 10
$ws.Range("AH28").Value = 67

# Row 31
$ws.Range("H31").Value = 3.75
$ws.Range("I31").Value = 4.65
$ws.Range("L31").Value = 1.23
$ws.Range("M31").Value = 3.35
$ws.Range("N31").Value = 1.7
$ws.Range("O31").Value = 1.91
$ws.Range("R31").Value = 1.7
$ws.Range("S31").Value = 1.91
$ws.Range("T31").Value = 7.3
$ws.Range("U31").Value = 8
$ws.Range("W31").Value = 12.5
$ws.Range("X31").Value = 12.5
$ws.Range("Y31").Value = 24
$ws.Range("Z31").Value = 11.75
$ws.Range("AA31").Value = 7.4
$ws.Range("AB31").Value = 15
$ws.Range("AC31").Value = 65
$ws.Range("AD31").Value = 450
$ws.Range("AE31").Value = 14.5
$ws.Range("AF31").Value = 29
$ws.Range("AJ31").Value = 45

# Row 32
$ws.Range("H32").Value = 4.05
$ws.Range("I32").Value = 4.55
$ws.Range("L32").Value = 1.21
$ws.Range("M32").Value = 3.9
$ws.Range("O32").Value = 2.12
$ws.Range("P32").Value = 1.31
$ws.Range("Q32").Value = 3.15
$ws.Range("R32").Value = 1.7
$ws.Range("S32").Value = 2.02
$ws.Range("T32").Value = 8.25
$ws.Range("W32").Value = 12.5
$ws.Range("AB32").Value = 15
$ws.Range("AC32").Value = 60
$ws.Range("AD32").Value = 400
$ws.Range("AE32").Value = 15
$ws.Range("AF32").Value = 27
$ws.Range("AG32").Value = 14.5

# Row 36
$ws.Range("H36").Value = 3.1
$ws.Range("L36").Value = 1.3
$ws.Range("M36").Value = 2.92
$ws.Range("N36").Value = 1.9
$ws.Range("O36").Value = 1.72
$ws.Range("P36").Value = 1.39
$ws.Range("Q36").Value = 2.55
$ws.Range("R36").Value = 1.7
$ws.Range("S36").Value = 1.93
$ws.Range("T36").Value = 8.5
$ws.Range("U36").Value = 13.5
$ws.Range("V36").Value = 9.75
$ws.Range("W36").Value = 30
$ws.Range("X36").Value = 22
$ws.Range("Y36").Value = 30
$ws.Range("Z36").Value = 9
$ws.Range("AB36").Value = 13
$ws.Range("AC36").Value = 60
$ws.Range("AD36").Value = 450
$ws.Range("AE36").Value = 8.25
$ws.Range("AF36").Value = 13
$ws.Range("AG36").Value = 9.5
$ws.Range("AH36").Value = 30
$ws.Range("AI36").Value = 22
$ws.Range("AJ36").Value = 30
